# copy (not move) loc data from tourney sheets. in prep for code change.
#
# Inserts two new rows right after the header row of the "Tournament" table
# (competition-key / host-key), and appends ten new rows at the bottom
# (venue-key.1 .. venue-key.10) holding the new "key" columns that mirror
# data that lives on the other sheets (Matches/Seeds/Colors). Only columns
# A and B are populated for all of the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tournament")
$lo = $ws.ListObjects.Item(1)

# insert two blank rows below the header, pushing all existing data rows
# down by two
$ws.Range("A2:A3").EntireRow.Insert()

# column A first for both new rows, then column B for both new rows
$ws.Range("A2").Value = "competition-key"
$ws.Range("A3").Value = "host-key"
$ws.Range("B2").Value = "womens-world-cup"
$ws.Range("B3").Value = "au-nz"

# append the ten venue-key rows at the bottom of the table
$ws.Range("A17").Value = "venue-key.1"
$ws.Range("B17").Value = "nz-auckland"

$ws.Range("A18").Value = "venue-key.2"
$ws.Range("B18").Value = "au-sydney_football"

$ws.Range("A19").Value = "venue-key.3"
$ws.Range("B19").Value = "nz-dunedin"

$ws.Range("A20").Value = "venue-key.4"
$ws.Range("B20").Value = "au-melbourne"

$ws.Range("A21").Value = "venue-key.5"
$ws.Range("B21").Value = "nz-wellington"

$ws.Range("A22").Value = "venue-key.6"
$ws.Range("B22").Value = "nz-hamilton"

$ws.Range("A23").Value = "venue-key.7"
$ws.Range("B23").Value = "au-brisbane"

$ws.Range("A24").Value = "venue-key.8"
$ws.Range("B24").Value = "au-perth"

$ws.Range("A25").Value = "venue-key.9"
$ws.Range("B25").Value = "au-adelaide"

$ws.Range("A26").Value = "venue-key.10"
$ws.Range("B26").Value = "au-sydney_australia"

# grow the table/autofilter to cover the new rows
$lo.Resize($ws.Range("A1:I26"))

# match the selection left behind in the saved file
$ws.Range("A2:XFD3").Select()
